$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 4741.6665
$ws.Range("I31").Value = 3290
$ws.Range("J31").Value = 12000
$ws.Range("K31").Value = 9870
$ws.Range("L31").Value = 36000
$ws.Range("M31").Value = -9640
$ws.Range("N31").Value = -36460
$ws.Range("H101").Value = 2625.7646
$ws.Range("I101").Value = 1094.4546
$ws.Range("J101").Value = 5433.1665
$ws.Range("K101").Value = 3283.3638
$ws.Range("L101").Value = 16299.4995
$ws.Range("M101").Value = -1661.3638
$ws.Range("N101").Value = -19543.4995
$ws.Range("H125").Value = 5506.4707
$ws.Range("I125").Value = 2078.375
$ws.Range("J125").Value = 8553.666999999999
$ws.Range("K125").Value = 18705.375
$ws.Range("L125").Value = 76983.003
$ws.Range("M125").Value = -16245.375
$ws.Range("N125").Value = -81903.003
$ws.Range("H129").Value = 2173.7637
$ws.Range("I129").Value = 1501.3636
$ws.Range("J129").Value = 3182.3635
$ws.Range("K129").Value = 4504.0908
$ws.Range("L129").Value = 9547.0905
$ws.Range("M129").Value = 495.9092000000001
$ws.Range("N129").Value = -19547.0905
$ws.Range("H132").Value = 158431.95
$ws.Range("I132").Value = 179249.64
$ws.Range("J132").Value = 9734.143
$ws.Range("K132").Value = 537748.92
$ws.Range("L132").Value = 29202.429
$ws.Range("M132").Value = -535218.92
$ws.Range("N132").Value = -34262.429
$ws.Range("H136").Value = 102094.75
$ws.Range("I136").Value = 80709
$ws.Range("J136").Value = 104038.91
$ws.Range("K136").Value = 80709
$ws.Range("L136").Value = 104038.91
$ws.Range("M136").Value = -75609
$ws.Range("N136").Value = -114238.91
$ws.Range("H137").Value = 297064.38
$ws.Range("I137").Value = 770965.25
$ws.Range("J137").Value = 3697.1904
$ws.Range("K137").Value = 2312895.75
$ws.Range("L137").Value = 11091.5712
$ws.Range("M137").Value = -2310345.75
$ws.Range("N137").Value = -16191.5712

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 31998.75
$ws.Range("I37").Value = 31998.75
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 31998.75
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -31725.75
$ws.Range("N37").Value = ""
$ws.Range("H61").Value = 5343.65
$ws.Range("I61").Value = 4940.4116
$ws.Range("J61").Value = 7628.6665
$ws.Range("K61").Value = 4940.4116
$ws.Range("L61").Value = 7628.6665
$ws.Range("M61").Value = -4728.4116
$ws.Range("N61").Value = -8052.6665
$ws.Range("H63").Value = 1000
$ws.Range("I63").Value = 1000
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1000
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -314
$ws.Range("H66").Value = 1000
$ws.Range("I66").Value = 1000
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 5000
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -1568
$ws.Range("H80").Value = 61998.332
$ws.Range("I80").Value = 50000
$ws.Range("J80").Value = 67997.5
$ws.Range("K80").Value = 50000
$ws.Range("L80").Value = 67997.5
$ws.Range("M80").Value = -49002
$ws.Range("N80").Value = -69993.5
$ws.Range("H83").Value = 61998.332
$ws.Range("I83").Value = 50000
$ws.Range("J83").Value = 67997.5
$ws.Range("K83").Value = 150000
$ws.Range("L83").Value = 203992.5
$ws.Range("M83").Value = -145008
$ws.Range("N83").Value = -213976.5
$ws.Range("H122").Value = 7090
$ws.Range("I122").Value = 4611.4546
$ws.Range("J122").Value = 10498
$ws.Range("K122").Value = 13834.3638
$ws.Range("L122").Value = 31494
$ws.Range("M122").Value = -11384.3638
$ws.Range("N122").Value = -36394
$ws.Range("H136").Value = 5343.65
$ws.Range("I136").Value = 4940.4116
$ws.Range("J136").Value = 7628.6665
$ws.Range("K136").Value = 14821.2348
$ws.Range("L136").Value = 22885.9995
$ws.Range("M136").Value = -12271.2348
$ws.Range("N136").Value = -27985.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 46500
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 46500
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 46500
$ws.Range("N35").Value = -47120
$ws.Range("H80").Value = 1275.6875
$ws.Range("I80").Value = 847.5
$ws.Range("J80").Value = 1418.4166
$ws.Range("K80").Value = 847.5
$ws.Range("L80").Value = 1418.4166
$ws.Range("M80").Value = 150.5
$ws.Range("N80").Value = -3414.4166
$ws.Range("H83").Value = 1275.6875
$ws.Range("I83").Value = 847.5
$ws.Range("J83").Value = 1418.4166
$ws.Range("K83").Value = 4237.5
$ws.Range("L83").Value = 7092.083000000001
$ws.Range("M83").Value = 754.5
$ws.Range("N83").Value = -17076.083
$ws.Range("H86").Value = 6013.5713
$ws.Range("I86").Value = 4200
$ws.Range("J86").Value = 6739
$ws.Range("K86").Value = 4200
$ws.Range("L86").Value = 6739
$ws.Range("M86").Value = -3077
$ws.Range("N86").Value = -8985
$ws.Range("H89").Value = 6013.5713
$ws.Range("I89").Value = 4200
$ws.Range("J89").Value = 6739
$ws.Range("K89").Value = 21000
$ws.Range("L89").Value = 33695
$ws.Range("M89").Value = -15384
$ws.Range("N89").Value = -44927
$ws.Range("H105").Value = 3618.2727
$ws.Range("I105").Value = 3410.6667
$ws.Range("J105").Value = 3696.125
$ws.Range("K105").Value = 3410.6667
$ws.Range("L105").Value = 3696.125
$ws.Range("M105").Value = -1663.6667
$ws.Range("N105").Value = -7190.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20002738
$ws.Range("I31").Value = 23257820
$ws.Range("J31").Value = 7234.7144
$ws.Range("K31").Value = 23257820
$ws.Range("L31").Value = 7234.7144
$ws.Range("M31").Value = -23257525
$ws.Range("N31").Value = -7824.7144
$ws.Range("H34").Value = 20002738
$ws.Range("I34").Value = 23257820
$ws.Range("J34").Value = 7234.7144
$ws.Range("K34").Value = 23257820
$ws.Range("L34").Value = 7234.7144
$ws.Range("M34").Value = -23257618
$ws.Range("N34").Value = -7638.7144
$ws.Range("H58").Value = 1724.6
$ws.Range("I58").Value = 1724.6
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1724.6
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -1521.6
$ws.Range("H94").Value = 1567.8096
$ws.Range("I94").Value = 438.4
$ws.Range("J94").Value = 1920.75
$ws.Range("K94").Value = 438.4
$ws.Range("L94").Value = 1920.75
$ws.Range("M94").Value = 12.60000000000002
$ws.Range("N94").Value = -2822.75
$ws.Range("H136").Value = 1724.6
$ws.Range("I136").Value = 1724.6
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5173.799999999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2623.799999999999
$ws.Range("H141").Value = 356167.84
$ws.Range("I141").Value = 99999.5
$ws.Range("J141").Value = 376661.3
$ws.Range("K141").Value = 99999.5
$ws.Range("L141").Value = 376661.3
$ws.Range("M141").Value = -94819.5
$ws.Range("N141").Value = -387021.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1635.9231
$ws.Range("I2").Value = 20.571428
$ws.Range("J2").Value = 3520.5
$ws.Range("K2").Value = 123.428568
$ws.Range("L2").Value = 21123
$ws.Range("M2").Value = -10.42856800000001
$ws.Range("N2").Value = -21349
$ws.Range("H4").Value = 4050948.2
$ws.Range("I4").Value = 158755.83
$ws.Range("J4").Value = 11835333
$ws.Range("K4").Value = 476267.49
$ws.Range("L4").Value = 35505999
$ws.Range("M4").Value = -476155.49
$ws.Range("N4").Value = -35506223
$ws.Range("H113").Value = 1619.84
$ws.Range("I113").Value = 597.25
$ws.Range("J113").Value = 1814.619
$ws.Range("K113").Value = 1791.75
$ws.Range("L113").Value = 5443.857
$ws.Range("M113").Value = 378.25
$ws.Range("N113").Value = -9783.857

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 9003223
$ws.Range("I11").Value = 3504833.8
$ws.Range("J11").Value = 20000000
$ws.Range("K11").Value = 3504833.8
$ws.Range("L11").Value = 20000000
$ws.Range("M11").Value = -3504694.8
$ws.Range("N11").Value = -20000278
$ws.Range("H46").Value = 20866.666
$ws.Range("I46").Value = 5500
$ws.Range("J46").Value = 23230.77
$ws.Range("K46").Value = 5500
$ws.Range("L46").Value = 23230.77
$ws.Range("M46").Value = -5344
$ws.Range("N46").Value = -23542.77
$ws.Range("H140").Value = 84998.57000000001
$ws.Range("I140").Value = 20000
$ws.Range("J140").Value = 95831.664
$ws.Range("K140").Value = 20000
$ws.Range("L140").Value = 95831.664
$ws.Range("M140").Value = -14820
$ws.Range("N140").Value = -106191.664
$ws.Range("H141").Value = 85000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 85000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 85000
$ws.Range("N141").Value = -95360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5625.067
$ws.Range("I7").Value = 4994.8
$ws.Range("J7").Value = 5940.2
$ws.Range("K7").Value = 4994.8
$ws.Range("L7").Value = 5940.2
$ws.Range("M7").Value = -4882.8
$ws.Range("N7").Value = -6164.2
$ws.Range("H55").Value = 616
$ws.Range("I55").Value = 288.16666
$ws.Range("J55").Value = 918.61536
$ws.Range("K55").Value = 288.16666
$ws.Range("L55").Value = 918.61536
$ws.Range("M55").Value = -115.16666
$ws.Range("N55").Value = -1264.61536
$ws.Range("H122").Value = 12498.167
$ws.Range("I122").Value = 5505.1665
$ws.Range("J122").Value = 19491.166
$ws.Range("K122").Value = 16515.4995
$ws.Range("L122").Value = 58473.49800000001
$ws.Range("M122").Value = -14065.4995
$ws.Range("N122").Value = -63373.49800000001
$ws.Range("H126").Value = 5625.067
$ws.Range("I126").Value = 4994.8
$ws.Range("J126").Value = 5940.2
$ws.Range("K126").Value = 14984.4
$ws.Range("L126").Value = 17820.6
$ws.Range("M126").Value = -12514.4
$ws.Range("N126").Value = -22760.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 22047.363
$ws.Range("I81").Value = 13421.75
$ws.Range("J81").Value = 26976.285
$ws.Range("K81").Value = 26843.5
$ws.Range("L81").Value = 53952.57
$ws.Range("M81").Value = -25782.5
$ws.Range("N81").Value = -56074.57
$ws.Range("H84").Value = 22047.363
$ws.Range("I84").Value = 13421.75
$ws.Range("J84").Value = 26976.285
$ws.Range("K84").Value = 134217.5
$ws.Range("L84").Value = 269762.85
$ws.Range("M84").Value = -128913.5
$ws.Range("N84").Value = -280370.85
$ws.Range("H107").Value = 928.05884
$ws.Range("I107").Value = 777.0714
$ws.Range("J107").Value = 1632.6666
$ws.Range("K107").Value = 2331.2142
$ws.Range("L107").Value = 4897.9998
$ws.Range("M107").Value = -411.2142000000003
$ws.Range("N107").Value = -8737.9998
$ws.Range("H126").Value = 111112390
$ws.Range("I126").Value = 125001310
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 375003930
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -375001460
$ws.Range("N126").Value = -7940
$ws.Range("H136").Value = 4891.239
$ws.Range("I136").Value = 2845.5745
$ws.Range("J136").Value = 9698.549999999999
$ws.Range("K136").Value = 8536.7235
$ws.Range("L136").Value = 29095.65
$ws.Range("M136").Value = -5986.7235
$ws.Range("N136").Value = -34195.64999999999
$ws.Range("H140").Value = 104999
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 104999
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 104999
$ws.Range("N140").Value = -115359
